# Updated cryptos list values (Price / Volume(1h)) per target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.195.53'
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").Value = '2.490.12'
$ws.Range("E3").Value = '  -1.62%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.05'
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.29'
$ws.Range("E6").Value = '  -7.01%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.517'
$ws.Range("E8").Value = '  -2.87%  '
$ws.Range("D9").Value = '2.488.79'
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("E10").Value = '  -4.03%  '
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.95'
$ws.Range("E13").Value = '  -3.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.95'
$ws.Range("E14").Value = '  -3.46%  '
$ws.Range("D15").Value = '2.941.15'
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("E16").Value = '  -3.58%  '
$ws.Range("D17").Value = '67.072.35'
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("D18").Value = '2.472.36'
$ws.Range("E18").Value = '  -1.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.63'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.83'
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '359.83'
$ws.Range("E21").Value = '  -2.32%  '
$ws.Range("E22").Value = '  -2.17%  '
$ws.Range("E23").Value = '  -6.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.54'
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("E26").Value = '  -4.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.43'
$ws.Range("E27").Value = '  -7.88%  '
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").Value = '2.614.47'
$ws.Range("E29").Value = '  -1.77%  '
$ws.Range("E30").Value = '  -6.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.05'
$ws.Range("E31").Value = '  -2.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '499.39'
$ws.Range("E32").Value = '  -7.77%  '
$ws.Range("E33").Value = '  -2.19%  '
$ws.Range("E34").Value = '  -5.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.127'
$ws.Range("E36").Value = '  -2.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.79'
$ws.Range("E37").Value = '  +1.45%  '
$ws.Range("E38").Value = '  +0.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.41'
$ws.Range("E39").Value = '  -3.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.56'
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("E41").Value = '  -4.32%  '
$ws.Range("E42").Value = '  -5.39%  '
$ws.Range("E43").Value = '  -5.70%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("E45").Value = '  -4.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.34'
$ws.Range("E46").Value = '  -1.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.90'
$ws.Range("E47").Value = '  -3.83%  '
$ws.Range("E48").Value = '  -2.84%  '
$ws.Range("E49").Value = '  -4.27%  '
$ws.Range("E50").Value = '  -5.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.65'
$ws.Range("E51").Value = '  -3.53%  '
